$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 keeps its original "text" representation ("001" -> "004"); force text
# format so Excel does not coerce it to the number 4, then restore the
# cell's default style so we don't leave a stray number-format override.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").Style = "Normal"

# Report date moves from 2018-12-31 to 2019-09-30 (kept as literal text,
# matching the source data's inline-string style for this column).
$ws.Range("N2").Value = "2019-09-30 00:00:00"

# Updated cash-flow figures / ratios for the new reporting period.
$ws.Range("O2").Value = 165540920.04
$ws.Range("P2").Value = 222.7763995116
$ws.Range("Q2").Value = 1130931107.76
$ws.Range("R2").Value = 1521.9485322512
$ws.Range("S2").Value = 40029174.83
$ws.Range("T2").Value = 53.8691910247
$ws.Range("U2").Value = -154231467.15
$ws.Range("V2").Value = -207.5567233453
$ws.Range("Y2").Value = 44077184.82
$ws.Range("Z2").Value = 59.3167932885
$ws.Range("AA2").Value = -84128592.59999999
$ws.Range("AB2").Value = -113.2159042663
$ws.Range("AC2").Value = -74308104.63

# AD2 (CCE_ADD_RATIO) no longer has a value for this period; clear it back
# to an empty cell.
$ws.Range("AD2").ClearContents()
